# Weekly update: insert two new price rows at the top of the existing
# "Comercializadora del Agro de Limarí - Uva" data block (rows 94-113),
# pushing the existing rows down to 96-115.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 94; this shifts old rows 94:113 down
# to 96:115 and extends the used range/dimension to A1:T115.
$ws.Rows("94:95").Insert()

# --- Row 94: new entry ---
$ws.Cells.Item(94, 1).Value = 2
$ws.Cells.Item(94, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(94, 3).Value = "Coquimbo"
$ws.Cells.Item(94, 4).Value2 = 44588
$ws.Cells.Item(94, 5).Value = 4
$ws.Cells.Item(94, 6).Value = "Fruta"
$ws.Cells.Item(94, 7).Value = 100109
$ws.Cells.Item(94, 8).Value = "Uva"
$ws.Cells.Item(94, 9).Value = 100109001
$ws.Cells.Item(94, 10).Value = "Uva"
$ws.Cells.Item(94, 11).Value = "Flame Seedless"
$ws.Cells.Item(94, 12).Value = "Primera"
$ws.Cells.Item(94, 13).Value = 700
$ws.Cells.Item(94, 14).Value = 4000
$ws.Cells.Item(94, 15).Value = 4500
$ws.Cells.Item(94, 16).Value = 4250
$ws.Cells.Item(94, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(94, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(94, 19).Value = 425
$ws.Cells.Item(94, 20).Value = 10

# --- Row 95: new entry ---
$ws.Cells.Item(95, 1).Value = 2
$ws.Cells.Item(95, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(95, 3).Value = "Coquimbo"
$ws.Cells.Item(95, 4).Value2 = 44588
$ws.Cells.Item(95, 5).Value = 4
$ws.Cells.Item(95, 6).Value = "Fruta"
$ws.Cells.Item(95, 7).Value = 100109
$ws.Cells.Item(95, 8).Value = "Uva"
$ws.Cells.Item(95, 9).Value = 100109001
$ws.Cells.Item(95, 10).Value = "Uva"
$ws.Cells.Item(95, 11).Value = "Superior Seedless"
$ws.Cells.Item(95, 12).Value = "Primera"
$ws.Cells.Item(95, 13).Value = 240
$ws.Cells.Item(95, 14).Value = 7000
$ws.Cells.Item(95, 15).Value = 7500
$ws.Cells.Item(95, 16).Value = 7250
$ws.Cells.Item(95, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(95, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(95, 19).Value = 725
$ws.Cells.Item(95, 20).Value = 10
